# TradingModel - 2021/11/16 data update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (A: id, B: Stock_Id, C: PositionSize) for rows 2..13
$data = @(
    @(1, 2314, 65),
    @(2, 2436, 60),
    @(4, 3035, 32),
    @(5, 3122, 90),
    @(6, 3141, 27),
    @(8, 3221, 61),
    @(10, 3588, 35),
    @(12, 6104, 36),
    @(13, 6138, 30),
    @(15, 6271, 20),
    @(16, 6411, 26),
    @(19, 8289, 160)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# Carry column A header style onto the newly added rows (11-13)
$ws.Range("A10").Copy()
$ws.Range("A11:A13").PasteSpecial(-4122)
